$d = $word.ActiveDocument

# --- Change 1: Replace meeting location text "LG 1 Canteen" -> "Library 1F 1-352" ---
# The replacement text needs a different run (larger font, sz=24 / 12pt) with NO
# character-width scaling (w:w), unlike the surrounding text which has w:w=105.
# Directly setting Font.Scaling on the target range is not effective in this
# runtime, so we stage the replacement text at a location that has clean
# (unscaled) formatting, capture its FormattedText, and transplant that onto
# the target range - this preserves the "no w:w" property correctly.

$docEnd = $d.Content.End
$stageRng = $d.Range($docEnd - 1, $docEnd - 1)
$stageRng.InsertAfter("Library 1F 1-352")
$stageRng.Font.Size = 12
$stagedFormattedText = $stageRng.FormattedText

$target = $d.Content
$target.Find.Execute("LG 1 Canteen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.FormattedText = $stagedFormattedText

# Remove the staging text we appended to the last paragraph, restoring it to "1/1"
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastParaRange = $lastPara.Range
$cleanupRng = $d.Range($lastParaRange.Start + 3, $lastParaRange.End - 1)
$cleanupRng.Delete()

# --- Change 2: Move the "_GoBack" bookmark from the end of the document to sit
# between "start " and "to " in the "APIs between those tasks..." paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$apiRng = $d.Content
$apiRng.Find.Execute("really start to ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $apiRng.Start + "really start ".Length
$bookmarkRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRng)
